$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13, shifting existing rows 13..109 down to 14..110.
$ws.Rows.Item(13).Insert()

# Populate the newly inserted row 13 with the new weekly price entry.
$ws.Range("A13").Value = 4
$ws.Range("B13").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C13").Value = "Los Lagos"
$ws.Range("D13").Value = 44473
$ws.Range("E13").Value = 10
$ws.Range("F13").Value = 100112009
$ws.Range("G13").Value = "Acelga"
$ws.Range("H13").Value = "Sin especificar"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 100
$ws.Range("K13").Value = 4000
$ws.Range("L13").Value = 4000
$ws.Range("M13").Value = 4000
$ws.Range("N13").Value = "$/docena de atados (4 kilos)"
$ws.Range("O13").Value = "Región del Maule"
$ws.Range("P13").Value = 1000
$ws.Range("Q13").Value = 4
$ws.Range("R13").Value = "Hortaliza"
